# 小组计划实施表 - add the 2017.11.22 (第十三周 周三) weekly block,
# and fill in the completion status for the previous (2017.11.20&11.21) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in "完成情况" (completion status) for the previous week's rows
#    (row 244-249), which had been left blank until now.
# ---------------------------------------------------------------------------
$ws.Range("C244").Value = '未完成'
$ws.Range("C245").Value = '未完成'
$ws.Range("C246").Value = '完成'
$ws.Range("C247").Value = '未完成'
$ws.Range("C248").Value = '完成'
$ws.Range("C249").Value = '未完成'

# ---------------------------------------------------------------------------
# 2) Append a brand-new weekly block starting at row 252 (row 251 stays
#    blank, matching the one-row gap used between every other weekly
#    block in this sheet).
#
#    Row 252         : date header line (merged A:D)              -> style like row 242
#    Row 253         : column headers (人员/计划任务/完成情况/备注) -> style like row 243
#    Rows 254-259     : the six team members' plan rows             -> style like rows 244-249
#    Row 260         : trailing "总结：" placeholder (merged A:D)   -> style like row 250
# ---------------------------------------------------------------------------

# Copy the row-format "templates" from the previous block (242-250) onto
# the new block (252-260) first, cell range by cell range, so the new
# rows end up referencing the very same cellXfs entries as their
# templates (s="2", s="3", s="12"/"5"/"5", etc.).
$ws.Range("A242:D242").Copy()
$ws.Range("A252:D252").PasteSpecial(-4122)

$ws.Range("A243:D243").Copy()
$ws.Range("A253:D253").PasteSpecial(-4122)

$ws.Range("A244:D244").Copy()
$ws.Range("A254:D254").PasteSpecial(-4122)

$ws.Range("A245:D245").Copy()
$ws.Range("A255:D255").PasteSpecial(-4122)

$ws.Range("A246:D246").Copy()
$ws.Range("A256:D256").PasteSpecial(-4122)

$ws.Range("A247:D247").Copy()
$ws.Range("A257:D257").PasteSpecial(-4122)

$ws.Range("A248:D248").Copy()
$ws.Range("A258:D258").PasteSpecial(-4122)

$ws.Range("A249:D249").Copy()
$ws.Range("A259:D259").PasteSpecial(-4122)

$ws.Range("A250:D250").Copy()
$ws.Range("A260:D260").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row heights, mirroring the source rows they were copied from.
$ws.Rows.Item(252).RowHeight = 22.5
$ws.Rows.Item(253).RowHeight = 22.5
$ws.Rows.Item(254).RowHeight = 45
$ws.Rows.Item(255).RowHeight = 22.5
$ws.Rows.Item(256).RowHeight = 45
$ws.Rows.Item(257).RowHeight = 22.5
$ws.Rows.Item(258).RowHeight = 45
$ws.Rows.Item(259).RowHeight = 45
$ws.Rows.Item(260).RowHeight = 22.5

# Merge the single-column banner rows (date header + trailing summary),
# same as every other weekly block in the sheet.
$ws.Range("A252:D252").Merge()
$ws.Range("A260:D260").Merge()

# --- Cell values -----------------------------------------------------------

# Row 252: new date header.
$ws.Range("A252").Value = '日期：2017.11.22 第十三周 周三'

# Row 253: column headers.
$ws.Range("A253").Value = '人员'
$ws.Range("B253").Value = '计划任务'
$ws.Range("C253").Value = '完成情况'
$ws.Range("D253").Value = '备注'

# Row 254: 李杰
$ws.Range("A254").Value = '李杰'
$ws.Range("B254").Value = '继续完善web app接口所有模块'

# Row 255: 周振朋
$ws.Range("A255").Value = '周振朋'
$ws.Range("B255").Value = '开发“买卖”模块'

# Row 256: 禤锦辉
$ws.Range("A256").Value = '禤锦辉'
$ws.Range("B256").Value = '帮助前端开发人员开发其中一个小模块'

# Row 257: 柯新钿
$ws.Range("A257").Value = '柯新钿'
$ws.Range("B257").Value = '完善前端“我的”模块'

# Row 258: 冯文雄
$ws.Range("A258").Value = '冯文雄'
$ws.Range("B258").Value = '帮助前端开发人员"消息"模块'

# Row 259: 阿卜力孜
$ws.Range("A259").Value = '阿卜力孜'
$ws.Range("B259").Value = '帮助前端开发人员开发其中一个小模块'

# Row 260: trailing summary placeholder for the new block.
$ws.Range("A260").Value = '总结：'

# ---------------------------------------------------------------------------
# 3) Leave the cursor where the author left off: cell C258.
# ---------------------------------------------------------------------------
$ws.Range("C258").Select()
